$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H54").Value = 17833.334
$ws.Range("I54").Value = 17833.334
$ws.Range("K54").Value = 17833.334
$ws.Range("M54").Value = -17347.334

$ws.Range("H129").Value = 1280.4117
$ws.Range("I129").Value = 532.1429000000001
$ws.Range("K129").Value = 1596.4287
$ws.Range("M129").Value = 3403.5713

$ws.Range("H132").Value = 1563.5555
$ws.Range("I132").Value = 1582
$ws.Range("K132").Value = 4746
$ws.Range("M132").Value = -2216

$ws.Range("H135").Value = 317.08
$ws.Range("I135").Value = 317.08
$ws.Range("K135").Value = 2853.72
$ws.Range("M135").Value = -318.7199999999998

$ws.Range("H137").Value = 3295.4226
$ws.Range("I137").Value = 2533.353
$ws.Range("J137").Value = 5238.7
$ws.Range("K137").Value = 7600.059
$ws.Range("L137").Value = 15716.1
$ws.Range("M137").Value = -5050.059
$ws.Range("N137").Value = -20816.1

$ws.Range("H141").Value = 30205.258
$ws.Range("I141").Value = 32641.322
$ws.Range("J141").Value = 11325.75
$ws.Range("K141").Value = 97923.966
$ws.Range("L141").Value = 33977.25
$ws.Range("M141").Value = -92743.966
$ws.Range("N141").Value = -44337.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4597.15
$ws.Range("I2").Value = 4744.1177
$ws.Range("K2").Value = 4744.1177
$ws.Range("M2").Value = -4631.1177

$ws.Range("H32").Value = 23266.588
$ws.Range("J32").Value = 121109.6
$ws.Range("L32").Value = 121109.6
$ws.Range("N32").Value = -121683.6

$ws.Range("H45").Value = 461350.3
$ws.Range("I45").Value = 778685.0600000001
$ws.Range("J45").Value = 2977.889
$ws.Range("K45").Value = 778685.0600000001
$ws.Range("L45").Value = 2977.889
$ws.Range("M45").Value = -778308.0600000001
$ws.Range("N45").Value = -3731.889

$ws.Range("H74").Value = 1291.4062
$ws.Range("J74").Value = 2037.1428
$ws.Range("L74").Value = 2037.1428
$ws.Range("N74").Value = -3785.1428

$ws.Range("H77").Value = 1291.4062
$ws.Range("J77").Value = 2037.1428
$ws.Range("L77").Value = 10185.714
$ws.Range("N77").Value = -18921.714

$ws.Range("H116").Value = 4597.15
$ws.Range("I116").Value = 4744.1177
$ws.Range("K116").Value = 4744.1177
$ws.Range("M116").Value = -2450.1177

$ws.Range("H132").Value = 1075.9062
$ws.Range("I132").Value = 1043.96
$ws.Range("J132").Value = 1190
$ws.Range("K132").Value = 3131.88
$ws.Range("L132").Value = 3570
$ws.Range("M132").Value = -601.8800000000001
$ws.Range("N132").Value = -8630

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4597.15
$ws.Range("I3").Value = 4744.1177
$ws.Range("K3").Value = 4744.1177
$ws.Range("M3").Value = -4630.1177

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 12226222
$ws.Range("J4").Value = 11670167
$ws.Range("L4").Value = 11670167
$ws.Range("N4").Value = -11670391

$ws.Range("H31").Value = 2261.5134
$ws.Range("I31").Value = 1419.6586
$ws.Range("K31").Value = 1419.6586
$ws.Range("M31").Value = -1124.6586

$ws.Range("H34").Value = 2261.5134
$ws.Range("I34").Value = 1419.6586
$ws.Range("K34").Value = 1419.6586
$ws.Range("M34").Value = -1217.6586

$ws.Range("H107").Value = 720.48
$ws.Range("I107").Value = 646.13794
$ws.Range("J107").Value = 823.1429000000001
$ws.Range("K107").Value = 646.13794
$ws.Range("L107").Value = 823.1429000000001
$ws.Range("M107").Value = 1273.86206
$ws.Range("N107").Value = -4663.1429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 750
$ws.Range("I80").Value = 750
$ws.Range("K80").Value = 2250
$ws.Range("M80").Value = -1314

$ws.Range("H83").Value = 750
$ws.Range("I83").Value = 750
$ws.Range("K83").Value = 6750
$ws.Range("M83").Value = -2070

$ws.Range("H129").Value = 1766.7858
$ws.Range("I129").Value = 490
$ws.Range("J129").Value = 3740
$ws.Range("K129").Value = 1470
$ws.Range("L129").Value = 11220
$ws.Range("M129").Value = 3530
$ws.Range("N129").Value = -21220

$ws.Range("H139").Value = 3317.842
$ws.Range("I139").Value = 1908.75
$ws.Range("J139").Value = 10833
$ws.Range("K139").Value = 5726.25
$ws.Range("L139").Value = 32499
$ws.Range("M139").Value = -586.25
$ws.Range("N139").Value = -42779

$ws.Range("H140").Value = 4023
$ws.Range("I140").Value = 3742.2222
$ws.Range("J140").Value = 4444.1665
$ws.Range("K140").Value = 11226.6666
$ws.Range("L140").Value = 13332.4995
$ws.Range("M140").Value = -6046.6666
$ws.Range("N140").Value = -23692.4995

$ws.Range("H141").Value = 2921
$ws.Range("I141").Value = 2839.2307
$ws.Range("J141").Value = 3984
$ws.Range("K141").Value = 8517.6921
$ws.Range("L141").Value = 11952
$ws.Range("M141").Value = -3337.6921
$ws.Range("N141").Value = -22312

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2385.3
$ws.Range("I102").Value = 2385.3
$ws.Range("K102").Value = 2385.3
$ws.Range("M102").Value = -763.3000000000002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 13015.346
$ws.Range("I20").Value = 2633.2778
$ws.Range("J20").Value = 36375
$ws.Range("K20").Value = 2633.2778
$ws.Range("L20").Value = 36375
$ws.Range("M20").Value = -2407.2778
$ws.Range("N20").Value = -36827

$ws.Range("H33").Value = 13159
$ws.Range("I33").Value = 8948.75
$ws.Range("K33").Value = 8948.75
$ws.Range("M33").Value = -8658.75

$ws.Range("H35").Value = 508.6
$ws.Range("I35").Value = 508.6
$ws.Range("K35").Value = 508.6
$ws.Range("M35").Value = -172.6

$ws.Range("H93").Value = 15470.833
$ws.Range("I93").Value = 1748.1666
$ws.Range("J93").Value = 56638.832
$ws.Range("K93").Value = 1748.1666
$ws.Range("L93").Value = 56638.832
$ws.Range("M93").Value = -500.1666
$ws.Range("N93").Value = -59134.832

$ws.Range("H132").Value = 3237.634
$ws.Range("I132").Value = 2781.7585
$ws.Range("J132").Value = 4339.3335
$ws.Range("K132").Value = 8345.2755
$ws.Range("L132").Value = 13018.0005
$ws.Range("M132").Value = -5815.2755
$ws.Range("N132").Value = -18078.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 55000
$ws.Range("J27").Value = 55000
$ws.Range("L27").Value = 55000
$ws.Range("M27").Value = -55138

$ws.Range("H100").Value = 2879.625
$ws.Range("I100").Value = 2738.2666
$ws.Range("K100").Value = 5476.5332
$ws.Range("M100").Value = -4935.5332

$ws.Range("I132").Value = 4325.4165
$ws.Range("K132").Value = 12976.2495
$ws.Range("M132").Value = -10446.2495

$ws.Range("H136").Value = 961.9583
$ws.Range("I136").Value = 574.15
$ws.Range("J136").Value = 2901
$ws.Range("K136").Value = 1722.45
$ws.Range("L136").Value = 8703
$ws.Range("M136").Value = 827.5500000000002
$ws.Range("N136").Value = -13803
